$wb = $excel.ActiveWorkbook

# Both the "展览" (Exhibitions) sheet and the "全部类型" (All types) sheet
# contain the same rows of data (rows 2-6) and both need to be updated
# to reflect the refreshed "想去人数" (want-to-go count) and
# "最低票价" (lowest ticket price) figures.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 6470
    $ws.Range("G2").Value = 55

    $ws.Range("F3").Value = 31

    $ws.Range("F4").Value = 188

    $ws.Range("F5").Value = 1019

    $ws.Range("F6").Value = 118
}

$wb.Save()
